$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2347.111
$ws.Range("I62").Value = 1341.8
$ws.Range("K62").Value = 1341.8
$ws.Range("M62").Value = -717.8
$ws.Range("H65").Value = 2347.111
$ws.Range("I65").Value = 1341.8
$ws.Range("K65").Value = 6709
$ws.Range("M65").Value = -3589
$ws.Range("H129").Value = 866.02563
$ws.Range("I129").Value = 347.16666
$ws.Range("J129").Value = 960.36365
$ws.Range("K129").Value = 1041.49998
$ws.Range("L129").Value = 2881.09095
$ws.Range("M129").Value = 3958.50002
$ws.Range("N129").Value = -12881.09095
$ws.Range("H132").Value = 141812.95
$ws.Range("I132").Value = 2689.3906
$ws.Range("K132").Value = 8068.1718
$ws.Range("M132").Value = -5538.1718
$ws.Range("H135").Value = 486
$ws.Range("I135").Value = 484.44446
$ws.Range("K135").Value = 4360.00014
$ws.Range("M135").Value = -1825.00014
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 13000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 3958.6028
$ws.Range("I32").Value = 3074.8406
$ws.Range("J32").Value = 19203.5
$ws.Range("K32").Value = 3074.8406
$ws.Range("L32").Value = 19203.5
$ws.Range("M32").Value = -2787.8406
$ws.Range("N32").Value = -19777.5
$ws.Range("H45").Value = 1647.4286
$ws.Range("I45").Value = 1647.4286
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1647.4286
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1270.4286
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 1178.4546
$ws.Range("I61").Value = 897.0769
$ws.Range("K61").Value = 897.0769
$ws.Range("M61").Value = -685.0769
$ws.Range("H132").Value = 1755.0145
$ws.Range("I132").Value = 1244.9615
$ws.Range("J132").Value = 3315.1765
$ws.Range("K132").Value = 3734.8845
$ws.Range("L132").Value = 9945.529500000001
$ws.Range("M132").Value = -1204.8845
$ws.Range("N132").Value = -15005.5295
$ws.Range("H136").Value = 1178.4546
$ws.Range("I136").Value = 897.0769
$ws.Range("K136").Value = 2691.2307
$ws.Range("M136").Value = -141.2307000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7694648
$ws.Range("I31").Value = 1210.5
$ws.Range("J31").Value = 16132612
$ws.Range("K31").Value = 1210.5
$ws.Range("L31").Value = 16132612
$ws.Range("M31").Value = -915.5
$ws.Range("N31").Value = -16133202
$ws.Range("H34").Value = 7694648
$ws.Range("I34").Value = 1210.5
$ws.Range("J34").Value = 16132612
$ws.Range("K34").Value = 1210.5
$ws.Range("L34").Value = 16132612
$ws.Range("M34").Value = -1008.5
$ws.Range("N34").Value = -16133016
$ws.Range("H58").Value = 1584.7079
$ws.Range("I58").Value = 1396.987
$ws.Range("J58").Value = 2789.25
$ws.Range("K58").Value = 1396.987
$ws.Range("L58").Value = 2789.25
$ws.Range("M58").Value = -1193.987
$ws.Range("N58").Value = -3195.25
$ws.Range("H94").Value = 1578.7059
$ws.Range("I94").Value = 927.75
$ws.Range("J94").Value = 1779
$ws.Range("K94").Value = 927.75
$ws.Range("L94").Value = 1779
$ws.Range("M94").Value = -476.75
$ws.Range("N94").Value = -2681
$ws.Range("H134").Value = 3125.5095
$ws.Range("I134").Value = 3259.3513
$ws.Range("J134").Value = 2816
$ws.Range("K134").Value = 9778.053899999999
$ws.Range("L134").Value = 8448
$ws.Range("M134").Value = -7243.053899999999
$ws.Range("N134").Value = -13518
$ws.Range("H136").Value = 1584.7079
$ws.Range("I136").Value = 1396.987
$ws.Range("J136").Value = 2789.25
$ws.Range("K136").Value = 4190.961
$ws.Range("L136").Value = 8367.75
$ws.Range("M136").Value = -1640.961
$ws.Range("N136").Value = -13467.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2412243
$ws.Range("I4").Value = 4822023.5
$ws.Range("J4").Value = 2462.6
$ws.Range("K4").Value = 14466070.5
$ws.Range("L4").Value = 7387.799999999999
$ws.Range("M4").Value = -14465958.5
$ws.Range("N4").Value = -7611.799999999999
$ws.Range("H113").Value = 629.86664
$ws.Range("I113").Value = 629.4545000000001
$ws.Range("K113").Value = 1888.3635
$ws.Range("M113").Value = 281.6364999999998
$ws.Range("H129").Value = 2467.0908
$ws.Range("I129").Value = 2930
$ws.Range("J129").Value = 2146.6155
$ws.Range("K129").Value = 8790
$ws.Range("L129").Value = 6439.8465
$ws.Range("M129").Value = -3790
$ws.Range("N129").Value = -16439.8465
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10999
$ws.Range("J5").Value = 10999
$ws.Range("L5").Value = 10999
$ws.Range("N5").Value = -11223
$ws.Range("H43").Value = 24266.715
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 24266.715
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 24266.715
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -24568.715
$ws.Range("H46").Value = 27420.875
$ws.Range("J46").Value = 30909.572
$ws.Range("L46").Value = 30909.572
$ws.Range("N46").Value = -31221.572
$ws.Range("H57").Value = 37247.5
$ws.Range("J57").Value = 36663.332
$ws.Range("L57").Value = 36663.332
$ws.Range("N57").Value = -38303.332
$ws.Range("H126").Value = 3220.47
$ws.Range("I126").Value = 2965.8691
$ws.Range("J126").Value = 4557.125
$ws.Range("K126").Value = 8897.6073
$ws.Range("L126").Value = 13671.375
$ws.Range("M126").Value = -6427.6073
$ws.Range("N126").Value = -18611.375
$ws.Range("H132").Value = 2092.68
$ws.Range("I132").Value = 1403.2285
$ws.Range("J132").Value = 3701.4
$ws.Range("K132").Value = 4209.6855
$ws.Range("L132").Value = 11104.2
$ws.Range("M132").Value = -1679.6855
$ws.Range("N132").Value = -16164.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2280.3096
$ws.Range("I136").Value = 1293.0294
$ws.Range("J136").Value = 6476.25
$ws.Range("K136").Value = 3879.0882
$ws.Range("L136").Value = 19428.75
$ws.Range("M136").Value = -1329.0882
$ws.Range("N136").Value = -24528.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 12000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 12000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 12000
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -12224
$ws.Range("H132").Value = 10754265
$ws.Range("I132").Value = 644.5833
$ws.Range("J132").Value = 47623820
$ws.Range("K132").Value = 1933.7499
$ws.Range("L132").Value = 142871460
$ws.Range("M132").Value = 596.2501
$ws.Range("N132").Value = -142876520
$ws.Range("H136").Value = 1810.0176
$ws.Range("J136").Value = 4228.2104
$ws.Range("L136").Value = 12684.6312
$ws.Range("N136").Value = -17784.6312
